# metricas_sel.xlsx update
# Adds 6 new activity rows (VectorColumna, VectorFila, setFila/setColumna,
# Producto MatrizMath, MatrizCuadrada, MatrizIdentidad) to the "Metricas"
# sheet, shifting the TOTALES / RESUMEN blocks down, refreshes the pie
# chart's source range + cached values, and repositions the chart to
# follow the moved rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: two new detail rows, inserted above the existing
#        "Producto/MatrizCuadrada/MatrizIdentidad" rows (old rows 12-13)
#        so the SUM()/shared-formula ranges below grow to include them. ---
$ws.Rows("12:13").Insert()

# Fix up the formatting of the two freshly inserted (blank) rows so they
# match their neighbours (row 11 style for row 13, row 16 style - which
# still carries the original "row 14" look - for row 12).
$ws.Range("A11:J11").Copy()
$ws.Range("A12:J13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A16").Copy()
$ws.Range("A12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Application.CutCopyMode = $false

# Row heights for the two new rows should match the surrounding data rows.
$ws.Rows("12:13").RowHeight = $ws.Rows("11").RowHeight

# --- 2. Fill in the six new activity rows. ---
# Row 10: Operaciones algebráicas en VectorColumna
$ws.Range("A10").Value = "Operaciones algebráicas en VectorColumna"
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 0.013888888888888888
$ws.Range("E10").Value = 0.7430555555555555
$ws.Range("F10").Value = 0.7604166666666666
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0

# Row 11: Operaciones algebráicas en VectorFila
$ws.Range("A11").Value = "Operaciones algebráicas en VectorFila"
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 12
$ws.Range("D11").Value = 0.006944444444444444
$ws.Range("E11").Value = 0.7645833333333334
$ws.Range("F11").Value = 0.7694444444444444
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

# Row 14: Creación de la Clase MatrizCuadrada
$ws.Range("A14").Value = "Creación de la Clase MatrizCuadrada"
$ws.Range("B14").Value = 20
$ws.Range("C14").Value = 18
$ws.Range("D14").Value = 0.010416666666666666
$ws.Range("E14").Value = 0.8590277777777778
$ws.Range("F14").Value = 0.8701388888888889
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0

# Row 15: Creación de la Clase MatrizIdentidad
$ws.Range("A15").Value = "Creación de la Clase MatrizIdentidad"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 0.001388888888888889
$ws.Range("E15").Value = 0.8722222222222222
$ws.Range("F15").Value = 0.875
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0

# Row 13: Producto entre MatrizMath con VectorFila y VectorColumna
$ws.Range("A13").Value = "Producto entre MatrizMath con VectorFila y VectorColumna"
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 24
$ws.Range("D13").Value = 0.010416666666666666
$ws.Range("E13").Value = 0.80625
$ws.Range("F13").Value = 0.8145833333333333
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# Row 12: Metodos setFila y setColumna de la Clase MatrizMath
$ws.Range("A12").Value = "Metodos setFila y setColumna de la Clase MatrizMath"
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 0.006944444444444444
$ws.Range("E12").Value = 0.7819444444444444
$ws.Range("F12").Value = 0.7909722222222223
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 0.010416666666666666

# The two inserted rows lost their G/J formulas on insert - restore them
# (same computation used by every other detail row: G = F-E, J = G+I).
$ws.Range("G12").Formula = "=F12-E12"
$ws.Range("J12").Formula = "=G12+I12"
$ws.Range("G13").Formula = "=F13-E13"
$ws.Range("J13").Formula = "=G13+I13"

# --- 3. Cosmetic sheet tweaks that accompanied the data entry. ---
$ws.Columns("A").ColumnWidth = 54.5703125
$ws.Columns("F").ColumnWidth = 11.5703125
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A14").Select()

# --- 4. Refresh the pie chart: source range moved from row 24/25 to
#        row 26/27, and the cached percentages follow the new totals. ---
$co = $ws.ChartObjects().Item(1)
$series = $co.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(,Metricas!`$A`$26:`$A`$27,Metricas!`$C`$26:`$C`$27,1)"

# --- 5. Move the chart down so it keeps sitting just under the RESUMEN
#        block, which is now two rows further down. ---
$topRow = $ws.Range("D20")
$bottomRow = $ws.Range("D42")
$co.Top = $topRow.Top + 12
$co.Height = ($bottomRow.Top + 12) - $co.Top
